$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the names/dates/types for the new rows first (6-8) ---
$ws.Cells.Item(6,1).Value = "2021-03-30"
$ws.Cells.Item(6,2).Value = "ERICSSON Next Generation Talent"
$ws.Cells.Item(6,3).Value = "Job application"
$ws.Cells.Item(6,4).Value = "Company webpage"

$ws.Cells.Item(7,1).Value = "2021-03-30"
$ws.Cells.Item(7,2).Value = "SEB Junior Data engineer"
$ws.Cells.Item(7,3).Value = "Job application"
$ws.Cells.Item(7,4).Value = "LinkedIn"

$ws.Cells.Item(8,1).Value = "2021-04-03"
$ws.Cells.Item(8,2).Value = "ATEA Intelligent Automation och RPA-konsult"
$ws.Cells.Item(8,3).Value = "Job application"
$ws.Cells.Item(8,4).Value = "Arbetsförmedlingen"

# --- Go back and update the RESULT_1 column for all applications that replied ---
$ws.Cells.Item(2,5).Value = "Replied"
$ws.Cells.Item(3,5).Value = "Replied"
$ws.Cells.Item(4,5).Value = "Replied"
$ws.Cells.Item(5,5).Value = "Replied"
$ws.Cells.Item(6,5).Value = "Replied"

# --- Update RESULT_2/RESULT_3 for the applications that progressed further ---
$ws.Cells.Item(2,6).Value = "Logic and personality test"
$ws.Cells.Item(2,7).Value = "Rejected"

$ws.Cells.Item(3,2).Value = "FOI Biträdande analytiker"
$ws.Cells.Item(3,6).Value = "Rejected"

$ws.Cells.Item(4,2).Value = "SCANIA Developer engineer"
$ws.Cells.Item(4,6).Value = "Interview 1"
$ws.Cells.Item(4,7).Value = "Rejected"

$ws.Cells.Item(5,2).Value = "SOLITA Data academy"
$ws.Cells.Item(5,4).Value = "LinkedIn"
$ws.Cells.Item(5,6).Value = "Logic and personality test"

$ws.Cells.Item(6,6).Value = "Interview 1"

# --- Remaining new rows (7-10) finish filling RESULT_1, plus row 9/10 names, and row 6 RESULT_3 ---
$ws.Cells.Item(7,5).Value = "Waiting for reply"
$ws.Cells.Item(8,5).Value = "Waiting for reply"

$ws.Cells.Item(9,1).Value = "2021-04-04"
$ws.Cells.Item(9,2).Value = "TietoEVRY"
$ws.Cells.Item(9,3).Value = "Job application"
$ws.Cells.Item(9,4).Value = "LinkedIn"
$ws.Cells.Item(9,5).Value = "Waiting for reply"

$ws.Cells.Item(6,7).Value = "Reference check"

$ws.Cells.Item(10,1).Value = "2021-04-13"
$ws.Cells.Item(10,2).Value = "SCANIA Developer engineer 2"
$ws.Cells.Item(10,3).Value = "Job application"
$ws.Cells.Item(10,4).Value = "Company webpage"
$ws.Cells.Item(10,5).Value = "Waiting for reply"

# --- Columns are now a bit wider to fit the new, longer content ---
$ws.Columns.Item(2).ColumnWidth = 41.5
$ws.Columns.Item(7).ColumnWidth = 14.833333333333334

# --- Selection moved to I5 as last action before save ---
$ws.Range("I5").Select()
